$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 205; this shifts the existing rows 205-292
# down to 206-293 and extends the sheet dimension to A1:T293.
$ws.Rows(205).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A205").Value2 = 10
$ws.Range("B205").Value2 = "Vega Modelo de Temuco"
$ws.Range("C205").Value2 = "La Araucanía"
$ws.Range("D205").Value2 = 44510
$ws.Range("E205").Value2 = 9
$ws.Range("F205").Value2 = "Fruta"
$ws.Range("G205").Value2 = 100103
$ws.Range("H205").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I205").Value2 = 100103006
$ws.Range("J205").Value2 = "Nectarín"
$ws.Range("K205").Value2 = "Early Glo"
$ws.Range("L205").Value2 = "Primera"
$ws.Range("M205").Value2 = 6500
$ws.Range("N205").Value2 = 34000
$ws.Range("O205").Value2 = 34000
$ws.Range("P205").Value2 = 34000
$ws.Range("Q205").Value2 = "$/bandeja 18 kilos granel"
$ws.Range("R205").Value2 = "Provincia de Quillota"
$ws.Range("S205").Value2 = 1889
$ws.Range("T205").Value2 = 18
